$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before F, shifting EmpID (old F) to H
$ws.Range("F:G").Insert()

# New header cells F1, G1 (bold + vertical-centered)
$ws.Range("F1").Value = "company / Initiative utilization "
$ws.Range("G1").Value = "Utilization on self "
$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").VerticalAlignment = -4108

# Row 2 values
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = [DateTime]::FromOADate(42257)
$ws.Range("C2").Value = [DateTime]::FromOADate(42257)
$ws.Range("D2").Value = "Project"
$ws.Range("E2").Value = "Billing Utilization"
$ws.Range("F2").Value = "company / Initiative utilization "
$ws.Range("G2").Value = "Utilization on self "
$ws.Range("H2").Value = 42

$ws.Range("F2:G2").VerticalAlignment = -4108

# Column widths (closest achievable values given the host's width quantization)
$ws.Columns("F").ColumnWidth = 28.333333333333332
$ws.Columns("G").ColumnWidth = 18
$ws.Columns("H").ColumnWidth = 8.833333333333334

$ws.Range("D2").Select()
